$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fishery_summary_table")

# Add the "Total" row beneath the summary table: a label in B22 and a
# SUM formula in C22 totalling the four fishery rows above it (C18:C21).
$ws.Range("B22").Value = "Total"
$ws.Range("C22").Formula = "=SUM(C18:C21)"

# Move the active selection to C22 (matches the saved selection in the diff)
$ws.Range("C22").Select()

# Don't recalculate before saving (mirrors calcOnSave="0" in the workbook)
$excel.CalculateBeforeSave = $false
